$d = $word.ActiveDocument

# The document ends with a ListParagraph-styled bullet paragraph
# ("...beta(a,b) for a>b?") that carries the "_GoBack" bookmark
# (collapsed) right after its final run. We need to append a brand new
# sibling bullet paragraph ("Also Alice, ...") after it, with the
# "_GoBack" bookmark ending up collapsed at the end of that new
# paragraph.
#
# Directly calling Bookmarks.Add with a collapsed Range positioned at
# the end of the (new) last paragraph mis-resolves in this host, so
# instead we: grow the existing bookmark's (collapsed) range in place
# with the new text via InsertBefore (Word keeps a collapsed bookmark
# sitting right after text inserted this way), then split the
# paragraph in two at the boundary between the old and new text. That
# leaves the old paragraph untouched and produces a new paragraph
# (inheriting the same ListParagraph/numPr formatting) that ends with
# the relocated, still-collapsed "_GoBack" bookmark - exactly matching
# native Word behaviour for this edit.

$newText = "Also Alice, in the data: the next does under ED30 is 44- is that the maximum dose now to give to the patients or can we give up to 58? I.E. can we give doses not given in the data or do we have to use those doses?"

if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks("_GoBack")
    $bmRange = $bm.Range
    $splitPos = $bmRange.Start

    $bmRange.InsertBefore($newText)

    $breakRange = $d.Range($splitPos, $splitPos)
    $breakRange.InsertParagraphAfter()
} else {
    # Fallback (no "_GoBack" bookmark present): just append a new
    # ListParagraph-styled bullet paragraph after the current last one.
    $lastPara = $d.Paragraphs.Last
    $lastPara.Range.InsertParagraphAfter()
    $d.Paragraphs.Last.Range.Text = $newText
}
